$d = $word.ActiveDocument

# 1. Remove the eastAsia-hint paragraph mark run formatting on the very
#    first paragraph (title paragraph) - drop its <w:pPr><w:rPr> font hint.
$p1 = $d.Paragraphs(1)
$p1.Range.Font.NameFarEast = ""

# 2. Append an ideographic full stop "。" run right after the sentence that
#    ends with "...穩定供氣，影響部分生態" (the "三接" question paragraph).
$r1 = $d.Content.Find.Execute("同意＝三接停工、另尋氣源；不同意＝三接續建、穩定供氣，影響部分生態", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($d.Content.Find.Found) {
    $insPoint = $d.Content
    $insPoint.Collapse(0)
}

$rng = $d.Content
$found = $rng.Find.Execute("同意＝三接停工、另尋氣源；不同意＝三接續建、穩定供氣，影響部分生態", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("。")
    $rng.Font.NameAscii = "+mn-ea"
}

# 3. Append an ideographic full stop "。" run right after the sentence that
#    ends with "...依然不可超過國際標準" (the "萊豬" question paragraph).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("同意＝萊豬不可進口，影響國際貿易；不同意＝萊豬繼續進口，依然不可超過國際標準", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter("。")
}
